$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay as text, matching source formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "93.211.45"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.412.26"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "232.80"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").Value = "618.70"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D7").Value = "1.42"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D8").Value = "0.391"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "0.973"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("D11").Value = "3.409.62"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "43.10"
$ws.Range("E12").Value = "  +5.71%  "
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "6.26"
$ws.Range("E14").Value = "  +3.10%  "
$ws.Range("D15").Value = "93.032.77"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "4.064.54"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "0.0000246"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "8.19"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("D19").Value = "3.412.02"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "17.95"
$ws.Range("E20").Value = "  +5.42%  "
$ws.Range("D21").Value = "11.63"
$ws.Range("E21").Value = "  +5.30%  "
$ws.Range("D22").Value = "0.490"
$ws.Range("E22").Value = "  +6.83%  "
$ws.Range("D23").Value = "3.38"
$ws.Range("E23").Value = "  +6.69%  "
$ws.Range("D24").Value = "496.14"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "6.67"
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("D26").Value = "0.0000183"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "95.08"
$ws.Range("E27").Value = "  +5.31%  "
$ws.Range("D28").Value = "11.97"
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "11.27"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.137"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.71"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D34").Value = "0.174"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").Value = "0.546"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "28.89"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "557.79"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "7.45"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "1.40"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.150"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.896"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "23.67"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "3.68"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "1.69"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0411"
$ws.Range("E46").Value = "  +3.47%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "5.42"
$ws.Range("E47").Value = "  -2.62%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "53.05"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "8.05"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "3.07"
$ws.Range("E51").Value = "  -2.40%  "
